# Auto update Excel log
# Appends 12 new sensor-reading rows (rows 2-13) to the PIR, Humidity and
# Temperature sheets of the SeniorConnect master log. Each sheet currently
# only has its header row (A1:F1); after the update the used range grows
# to A1:F13.

$wb = $excel.ActiveWorkbook

# Columns: Date, Timestamp, Hour, Location, Value, Status
$pirData = @(
    @("2026-01-28","14:39:31","14:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","14:39:36","14:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","14:39:41","14:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","14:39:47","14:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","14:39:51","14:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","14:39:57","14:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","14:40:02","14:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","14:40:07","14:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","14:40:12","14:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","14:40:17","14:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","14:40:22","14:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","14:40:27","14:00","Bathroom","No Motion","Inactive")
)

$humidityData = @(
    @("2026-01-28","14:39:31","14:00","Bathroom","88.8%","Active"),
    @("2026-01-28","14:39:39","14:00","Bathroom","87.8%","Active"),
    @("2026-01-28","14:39:43","14:00","Bathroom","88.8%","Active"),
    @("2026-01-28","14:39:47","14:00","Bathroom","87.8%","Active"),
    @("2026-01-28","14:39:51","14:00","Bathroom","88.8%","Active"),
    @("2026-01-28","14:39:59","14:00","Bathroom","88.8%","Active"),
    @("2026-01-28","14:40:03","14:00","Bathroom","88.7%","Active"),
    @("2026-01-28","14:40:07","14:00","Bathroom","87.9%","Active"),
    @("2026-01-28","14:40:11","14:00","Bathroom","88.8%","Active"),
    @("2026-01-28","14:40:19","14:00","Bathroom","87.9%","Active"),
    @("2026-01-28","14:40:23","14:00","Bathroom","88.8%","Active"),
    @("2026-01-28","14:40:27","14:00","Bathroom","87.9%","Active")
)

$temperatureData = @(
    @("2026-01-28","14:39:31","14:00","Bathroom","22.7C","Active"),
    @("2026-01-28","14:39:39","14:00","Bathroom","22.7C","Active"),
    @("2026-01-28","14:39:43","14:00","Bathroom","22.7C","Active"),
    @("2026-01-28","14:39:47","14:00","Bathroom","22.7C","Active"),
    @("2026-01-28","14:39:51","14:00","Bathroom","22.7C","Active"),
    @("2026-01-28","14:39:59","14:00","Bathroom","22.7C","Active"),
    @("2026-01-28","14:40:03","14:00","Bathroom","22.7C","Active"),
    @("2026-01-28","14:40:07","14:00","Bathroom","22.7C","Active"),
    @("2026-01-28","14:40:11","14:00","Bathroom","22.7C","Active"),
    @("2026-01-28","14:40:19","14:00","Bathroom","22.7C","Active"),
    @("2026-01-28","14:40:23","14:00","Bathroom","22.7C","Active"),
    @("2026-01-28","14:40:27","14:00","Bathroom","22.7C","Active")
)

function Write-SheetRows {
    param($SheetName, $Rows)

    $ws = $wb.Worksheets.Item($SheetName)

    # Pre-format column A (dates) as Text so values like "2026-01-28" are
    # stored as literal strings instead of being auto-converted to date
    # serial numbers.
    $lastRow = 1 + $Rows.Count
    $ws.Range("A2:A$lastRow").NumberFormat = "@"

    # Humidity's Value column holds strings like "88.8%" which Excel would
    # otherwise auto-convert to a percentage number - keep it textual too.
    if ($SheetName -eq "Humidity") {
        $ws.Range("E2:E$lastRow").NumberFormat = "@"
    }

    for ($i = 0; $i -lt $Rows.Count; $i++) {
        $r = $i + 2
        $rowValues = $Rows[$i]
        for ($j = 0; $j -lt $rowValues.Count; $j++) {
            $c = $j + 1
            $ws.Cells.Item($r, $c).Value = $rowValues[$j]
        }
    }
}

Write-SheetRows "PIR" $pirData
Write-SheetRows "Humidity" $humidityData
Write-SheetRows "Temperature" $temperatureData
